$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend rows 58:64 with the same formatting (styles) as row 57,
# which already has the date-format (col A) and number-format (cols B:E) styles applied.
$ws.Range("A57:E57").Copy()
$ws.Range("A58:E64").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update values (ECB May 2025 data refresh).
$ws.Cells.Item(5, 3).Value2 = 0.77241760492324829
$ws.Cells.Item(5, 5).Value2 = 0.62814927101135254
$ws.Cells.Item(6, 3).Value2 = 0.35116145014762878
$ws.Cells.Item(6, 5).Value2 = 0.63088101148605347
$ws.Cells.Item(7, 3).Value2 = 0.56391721963882446
$ws.Cells.Item(7, 5).Value2 = 0.66285288333892822
$ws.Cells.Item(8, 3).Value2 = 0.51864200830459595
$ws.Cells.Item(8, 5).Value2 = 0.75899547338485718
$ws.Cells.Item(9, 3).Value2 = 0.93460792303085327
$ws.Cells.Item(9, 5).Value2 = 0.75947266817092896
$ws.Cells.Item(10, 3).Value2 = 0.64454001188278198
$ws.Cells.Item(10, 5).Value2 = 0.67708730697631836
$ws.Cells.Item(11, 3).Value2 = 0.85468411445617676
$ws.Cells.Item(11, 5).Value2 = 0.69506186246871948
$ws.Cells.Item(12, 3).Value2 = 1.4319933652877808
$ws.Cells.Item(12, 5).Value2 = 0.69977849721908569
$ws.Cells.Item(13, 3).Value2 = 0.7632904052734375
$ws.Cells.Item(13, 5).Value2 = 0.70903778076171875
$ws.Cells.Item(14, 3).Value2 = 0.030949430540204048
$ws.Cells.Item(14, 5).Value2 = 0.68135523796081543
$ws.Cells.Item(15, 3).Value2 = 0.5129324197769165
$ws.Cells.Item(15, 5).Value2 = 0.67107510566711426
$ws.Cells.Item(16, 3).Value2 = 0.60636657476425171
$ws.Cells.Item(16, 5).Value2 = 0.6263502836227417
$ws.Cells.Item(17, 3).Value2 = 0.60197561979293823
$ws.Cells.Item(17, 5).Value2 = 0.54420250654220581
$ws.Cells.Item(18, 3).Value2 = 0.68546545505523682
$ws.Cells.Item(18, 5).Value2 = 0.50154620409011841
$ws.Cells.Item(19, 3).Value2 = 0.55201870203018188
$ws.Cells.Item(19, 5).Value2 = 0.5723302960395813
$ws.Cells.Item(20, 3).Value2 = 0.45216038823127747
$ws.Cells.Item(20, 5).Value2 = 0.57996618747711182
$ws.Cells.Item(21, 3).Value2 = 0.69266366958618164
$ws.Cells.Item(21, 5).Value2 = 0.55501270294189453
$ws.Cells.Item(22, 3).Value2 = 0.37938377261161804
$ws.Cells.Item(22, 5).Value2 = 0.48392975330352783
$ws.Cells.Item(23, 3).Value2 = 0.66800594329833984
$ws.Cells.Item(23, 5).Value2 = 0.42876240611076355
$ws.Cells.Item(24, 3).Value2 = 0.58165568113327026
$ws.Cells.Item(24, 5).Value2 = 0.35514086484909058
$ws.Cells.Item(25, 3).Value2 = 0.38178521394729614
$ws.Cells.Item(25, 4).Value2 = 0.38858711719512939
$ws.Cells.Item(25, 5).Value2 = 0.29685455560684204
$ws.Cells.Item(26, 3).Value2 = -0.037770979106426239
$ws.Cells.Item(26, 5).Value2 = 0.19572587311267853
$ws.Cells.Item(27, 3).Value2 = 0.18895919620990753
$ws.Cells.Item(27, 5).Value2 = 0.16923424601554871
$ws.Cells.Item(28, 3).Value2 = -0.11057507991790771
$ws.Cells.Item(28, 5).Value2 = 0.084503613412380219
$ws.Cells.Item(29, 3).Value2 = -0.072416514158248901
$ws.Cells.Item(29, 5).Value2 = 0.0065942686051130295
$ws.Cells.Item(30, 3).Value2 = -0.21749438345432281
$ws.Cells.Item(30, 5).Value2 = -0.070078320801258087
$ws.Cells.Item(31, 3).Value2 = 0.14095918834209442
$ws.Cells.Item(31, 5).Value2 = -0.10176176577806473
$ws.Cells.Item(32, 3).Value2 = -0.09456980973482132
$ws.Cells.Item(32, 5).Value2 = -0.14302507042884827
$ws.Cells.Item(33, 3).Value2 = -0.11952841281890869
$ws.Cells.Item(33, 5).Value2 = -0.15402337908744812
$ws.Cells.Item(34, 3).Value2 = -0.30826810002326965
$ws.Cells.Item(34, 5).Value2 = -0.19229646027088165
$ws.Cells.Item(35, 3).Value2 = -0.32292196154594421
$ws.Cells.Item(35, 5).Value2 = -0.17741735279560089
$ws.Cells.Item(36, 3).Value2 = -0.18241055309772491
$ws.Cells.Item(36, 5).Value2 = -0.2225511372089386
$ws.Cells.Item(37, 3).Value2 = -0.20955987274646759
$ws.Cells.Item(37, 5).Value2 = -0.21012218296527863
$ws.Cells.Item(38, 3).Value2 = -0.41687420010566711
$ws.Cells.Item(38, 5).Value2 = -0.218982994556427
$ws.Cells.Item(39, 3).Value2 = -0.083582490682601929
$ws.Cells.Item(39, 5).Value2 = -0.18898722529411316
$ws.Cells.Item(40, 3).Value2 = -0.26524484157562256
$ws.Cells.Item(40, 5).Value2 = -0.15824517607688904
$ws.Cells.Item(41, 3).Value2 = 0.017290746793150902
$ws.Cells.Item(41, 5).Value2 = -0.14717511832714081
$ws.Cells.Item(42, 3).Value2 = -0.19927568733692169
$ws.Cells.Item(42, 5).Value2 = -0.095496423542499542
$ws.Cells.Item(43, 3).Value2 = -0.03830612450838089
$ws.Cells.Item(43, 5).Value2 = -0.024876292794942856
$ws.Cells.Item(44, 3).Value2 = -0.046243507415056229
$ws.Cells.Item(44, 5).Value2 = 0.011214782483875751
$ws.Cells.Item(45, 3).Value2 = -0.082780137658119202
$ws.Cells.Item(45, 5).Value2 = 0.069724157452583313
$ws.Cells.Item(46, 3).Value2 = 0.25554844737052917
$ws.Cells.Item(46, 5).Value2 = 0.085541516542434692
$ws.Cells.Item(47, 3).Value2 = 0.21870696544647217
$ws.Cells.Item(47, 5).Value2 = 0.11117653548717499
$ws.Cells.Item(48, 3).Value2 = 0.24123717844486237
$ws.Cells.Item(48, 5).Value2 = 0.13925060629844666
$ws.Cells.Item(49, 3).Value2 = 0.26133951544761658
$ws.Cells.Item(49, 5).Value2 = 0.16183334589004517
$ws.Cells.Item(50, 3).Value2 = 0.15964698791503906
$ws.Cells.Item(50, 5).Value2 = 0.1998257040977478
$ws.Cells.Item(51, 3).Value2 = 0.031439471989870071
$ws.Cells.Item(51, 5).Value2 = 0.17734892666339874
$ws.Cells.Item(52, 3).Value2 = 0.21436057984828949
$ws.Cells.Item(52, 4).Value2 = 0.31449529528617859
$ws.Cells.Item(52, 5).Value2 = 0.18676285445690155
$ws.Cells.Item(53, 3).Value2 = 0.15700110793113708
$ws.Cells.Item(53, 4).Value2 = 0.34667885303497314
$ws.Cells.Item(53, 5).Value2 = 0.21741728484630585
$ws.Cells.Item(54, 3).Value2 = 0.25915110111236572
$ws.Cells.Item(54, 4).Value2 = 0.36161810159683228
$ws.Cells.Item(54, 5).Value2 = 0.2290860116481781
$ws.Cells.Item(55, 3).Value2 = 0.053257469087839127
$ws.Cells.Item(55, 4).Value2 = 0.37769380211830139
$ws.Cells.Item(55, 5).Value2 = 0.24751675128936768
$ws.Cells.Item(56, 2).Value2 = 0.4117148220539093
$ws.Cells.Item(56, 3).Value2 = 0.30343231558799744
$ws.Cells.Item(56, 4).Value2 = 0.44006979465484619
$ws.Cells.Item(56, 5).Value2 = 0.31197571754455566
$ws.Cells.Item(57, 2).Value2 = 0.64497452974319458
$ws.Cells.Item(57, 3).Value2 = 0.51712703704833984
$ws.Cells.Item(57, 4).Value2 = 0.45727953314781189
$ws.Cells.Item(57, 5).Value2 = 0.3360971212387085
$ws.Cells.Item(58, 1).Value2 = 45536
$ws.Cells.Item(58, 2).Value2 = 0.48986724019050598
$ws.Cells.Item(58, 3).Value2 = 0.36635807156562805
$ws.Cells.Item(58, 4).Value2 = 0.48500567674636841
$ws.Cells.Item(58, 5).Value2 = 0.36640718579292297
$ws.Cells.Item(59, 1).Value2 = 45566
$ws.Cells.Item(59, 2).Value2 = 0.44431072473526001
$ws.Cells.Item(59, 3).Value2 = 0.32552364468574524
$ws.Cells.Item(59, 4).Value2 = 0.48190978169441223
$ws.Cells.Item(59, 5).Value2 = 0.36044332385063171
$ws.Cells.Item(60, 1).Value2 = 45597
$ws.Cells.Item(60, 2).Value2 = 0.71741873025894165
$ws.Cells.Item(60, 3).Value2 = 0.61157023906707764
$ws.Cells.Item(60, 4).Value2 = 0.4954032301902771
$ws.Cells.Item(60, 5).Value2 = 0.37740015983581543
$ws.Cells.Item(61, 1).Value2 = 45627
$ws.Cells.Item(61, 2).Value2 = 0.54052573442459106
$ws.Cells.Item(61, 3).Value2 = 0.43145304918289185
$ws.Cells.Item(61, 4).Value2 = 0.50586426258087158
$ws.Cells.Item(61, 5).Value2 = 0.38664612174034119
$ws.Cells.Item(62, 1).Value2 = 45658
$ws.Cells.Item(62, 2).Value2 = 0.55239719152450562
$ws.Cells.Item(62, 3).Value2 = 0.42979174852371216
$ws.Cells.Item(62, 4).Value2 = 0.48599138855934143
$ws.Cells.Item(62, 5).Value2 = 0.36800599098205566
$ws.Cells.Item(63, 1).Value2 = 45689
$ws.Cells.Item(63, 2).Value2 = 0.34792107343673706
$ws.Cells.Item(63, 3).Value2 = 0.20547620952129364
$ws.Cells.Item(63, 4).Value2 = 0.48534542322158813
$ws.Cells.Item(63, 5).Value2 = 0.36828064918518066
$ws.Cells.Item(64, 1).Value2 = 45717
$ws.Cells.Item(64, 2).Value2 = 0.30949902534484863
$ws.Cells.Item(64, 3).Value2 = 0.20586903393268585
$ws.Cells.Item(64, 4).Value2 = 0.49355235695838928
$ws.Cells.Item(64, 5).Value2 = 0.37683206796646118
